$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1: bold font, thin box border, centered/top aligned
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Borders.LineStyle = 1
$b1.Borders.Weight = 2
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160

# Copy the exact same formatting onto A2 (using copy/paste-formats keeps
# the style table from growing extra orphaned entries that per-cell
# property assignment would otherwise introduce)
$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)
$excel.CutCopyMode = $false
